# executive-presentation.pptx edit
#
# Summary of target state (per the authoritative diff):
#   - The deck is trimmed from 17 slides down to 4 slides.
#   - Kept (by original content/position): slide 1 (title slide),
#     slide 7 ("Why This Solution?"), slide 8 ("Business Value -
#     Financial Impact"), slide 11 ("Risk Mitigation"). All other
#     slides (2-6, 9-10, 12-17) are removed.
#   - After trimming, those four kept slides naturally land at
#     positions 1, 2, 3, 4 (their relative order is preserved).
#   - On the new slide 2 (was "Why This Solution?"): the last
#     comparison-table row ("[Current limitation 3]" / "[Our
#     advantage 3]") loses its bold emphasis.
#   - On the new slide 3 (was "Business Value - Financial Impact"):
#     the header row ("Metric" / "Value") and the "ROI" /
#     "[Percentage]" row lose their bold emphasis.
#   - On the new slide 4 (was "Risk Mitigation"): the header row
#     ("Risk" / "Mitigation Strategy" / "Success Probability") and
#     the "[Risk 3]" row lose their bold emphasis.

$p = $ppt.ActivePresentation

# --- 1. Remove the unwanted slides -----------------------------------
# Delete from the highest index down to the lowest so earlier deletes
# never shift the index of a slide we still need to remove.
# Slides kept (by their current 1-based index): 1, 7, 8, 11.
$slidesToDelete = @(17,16,15,14,13,12,10,9,6,5,4,3,2)
foreach ($idx in $slidesToDelete) {
    $p.Slides.Item($idx).Delete()
}

# After the deletions above, the presentation has exactly 4 slides,
# in this order: [old slide 1, old slide 7, old slide 8, old slide 11]
# -> new positions 1, 2, 3, 4 respectively.

# --- helper: strip bold from a table cell's text ----------------------
function Clear-CellBold($table, $row, $col) {
    $table.Cell($row, $col).Shape.TextFrame.TextRange.Font.Bold = $false
}

# --- 2. New slide 2 (was "Why This Solution?") ------------------------
$s2 = $p.Slides.Item(2)
$tbl2 = $s2.Shapes.Item(3).Table
# Row 4 = "[Current limitation 3]" / "[Our advantage 3]"
Clear-CellBold $tbl2 4 1
Clear-CellBold $tbl2 4 2

# --- 3. New slide 3 (was "Business Value - Financial Impact") --------
$s3 = $p.Slides.Item(3)
$tbl3 = $s3.Shapes.Item(3).Table
# Row 1 = header ("Metric" / "Value")
Clear-CellBold $tbl3 1 1
Clear-CellBold $tbl3 1 2
# Row 6 = "ROI" / "[Percentage]"
Clear-CellBold $tbl3 6 1
Clear-CellBold $tbl3 6 2

# --- 4. New slide 4 (was "Risk Mitigation") ---------------------------
$s4 = $p.Slides.Item(4)
$tbl4 = $s4.Shapes.Item(3).Table
# Row 1 = header ("Risk" / "Mitigation Strategy" / "Success Probability")
Clear-CellBold $tbl4 1 1
Clear-CellBold $tbl4 1 2
Clear-CellBold $tbl4 1 3
# Row 4 = "[Risk 3]" / "[Mitigation approach]" / "[Percentage]"
Clear-CellBold $tbl4 4 1
Clear-CellBold $tbl4 4 2
Clear-CellBold $tbl4 4 3

Write-Output "Final slide count: $($p.Slides.Count)"
